$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D4: "Triadic closure, strength of weak ties" -> "Triadic closure"
$ws.Range("D4").Value = "Triadic closure"

# G4: add new homework link
$ws.Range("G4").Value = "[Hwk 2: Analyzing personal network data](https://datahub.berkeley.edu/hub/user-redirect/git-pull?repo=https%3A%2F%2Fgithub.com%2Fdfeehan%2Fdemog180-fa2024&branch=main&urlpath=tree%2Fdemog180-fa2024%2Fhwk%2Fhwk02%2Fhwk02.ipynb)"

# D5: "Social capital and structural holes" -> "Strength of Weak Ties, Social Capital, Structural Holes"
$ws.Range("D5").Value = "Strength of Weak Ties, Social Capital, Structural Holes"

# G6: remove the old homework entry (no link); it has been moved/replaced by G4
$ws.Range("G6").ClearContents()
